$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2698
$ws.Range("F3").Value = 1043
$ws.Range("F4").Value = 19297
$ws.Range("F6").Value = 2183
$ws.Range("F7").Value = 739
$ws.Range("F8").Value = 610
$ws.Range("F10").Value = 675
$ws.Range("F12").Value = 243
$ws.Range("F14").Value = 354
$ws.Range("F15").Value = 68
$ws.Range("F18").Value = 175
$ws.Range("F19").Value = 17
$ws.Range("F20").Value = 19
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 120
$ws.Range("F7").Value = 271
$ws.Range("F8").Value = 125
$ws.Range("F10").Value = 12
$ws.Range("F15").Value = 59
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5972
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 5972
$ws.Range("F7").Value = 2698
$ws.Range("F8").Value = 1043
$ws.Range("F9").Value = 19297
$ws.Range("F13").Value = 120
$ws.Range("F14").Value = 271
$ws.Range("F15").Value = 2183
$ws.Range("F16").Value = 739
$ws.Range("F17").Value = 125
$ws.Range("F18").Value = 610
$ws.Range("F20").Value = 675
$ws.Range("F22").Value = 243
$ws.Range("F25").Value = 12
$ws.Range("F27").Value = 354
$ws.Range("F28").Value = 68
$ws.Range("F34").Value = 175
$ws.Range("F35").Value = 59
$ws.Range("F37").Value = 17
$ws.Range("F40").Value = 19
